# Apply company-name corrections (remove trailing ",THE" and insert
# letter-spacing in certain acronym based bank names), then adjust the
# window/selection state to reflect the saved view from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = "KARUR VYSYA BANK LTD."
$ws.Range("A4").Value  = "I C I C I BANK LTD."
$ws.Range("A6").Value  = "LAKSHMI VILAS BANK LTD."
$ws.Range("A8").Value  = "D C B BANK LTD."
$ws.Range("A14").Value = "H D F C BANK LTD."
$ws.Range("A19").Value = "I D B I BANK LTD."
$ws.Range("A21").Value = "FEDERAL BANK LTD."
$ws.Range("A34").Value = "KARNATAKA BANK LTD."
$ws.Range("A40").Value = "SOUTH INDIAN BANK LTD."
$ws.Range("A41").Value = "JAMMU & KASHMIR BANK LTD."

# Restore the selection that was active when the author last saved
# the file.
$ws.Range("D30").Select()
